# Plantilla Lista de Tareas de la 5ta Iteracion
# Commit: "Mecanismo para obtener fotos de alumnos."
#
# Changes applied:
#  1. Casos de Uso!N10 gets a value of 2 (previously empty), which ripples
#     through the shared "restante"/"total" formulas for the rest of row 10.
#  2. Casos de Uso!D11 ("Mostrar unicamente a los alumnos activos para los
#     profesores.") is given an underlined font, to flag/highlight the new
#     task about obtaining student photos.
#  3. The frozen-pane scroll position / last selected cell on the sheet are
#     updated to reflect where the author ended up working (row 6 at the
#     top of the scrollable area, D11 selected).
#  4. The five right-most quarter-merges in the header row (row 4) get
#     touched again (un-merge/re-merge) which is what shuffles their order
#     to the end of the <mergeCells> list, matching the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# --- 1. Data edit: register 2 consumed hours for the task in row 10 on Day 3 (column N) ---
$ws.Range("N10").Value = 2

# --- 2. Formatting edit: underline the task description for row 11 ---
$ws.Range("D11").Font.Underline = $true

# --- 3. Touch the trailing header merges so they get re-appended (matches
#        the reordering seen in the saved workbook) ---
$ws.Range("AZ4:BA4").UnMerge()
$ws.Range("AO4:AP4").UnMerge()
$ws.Range("AR4:AS4").UnMerge()
$ws.Range("AU4:AV4").UnMerge()
$ws.Range("AX4:AY4").UnMerge()
$ws.Range("AZ4:BA4").Merge()
$ws.Range("AO4:AP4").Merge()
$ws.Range("AR4:AS4").Merge()
$ws.Range("AU4:AV4").Merge()
$ws.Range("AX4:AY4").Merge()

# --- 4. View state: scroll the frozen window up so row 6 is the first
#        visible row under the freeze, and finish with D11 selected ---
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 7
$ws.Range("D11").Select()
